# Active_Outages.xlsx update - 6/18/2025, 5:08:02 PM
#
# 1) Sheet "R1": remove the stale outage row (R5 / HAJ0155 / Good) that was
#    sitting at the bottom of the table (row 6), shrinking the used range
#    from A1:L6 to A1:L5.
# 2) Refresh the "Elapsed Duration(Hrs)" (column G) values across every
#    region sheet to reflect the passage of time since the report was
#    generated.

$wb = $excel.ActiveWorkbook

# --- R1: drop the resolved HAJ0155 row, then bump the Elapsed Duration values ---
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Rows.Item(6).Delete()
$ws1.Range("G2").Value = "3930:22:13"
$ws1.Range("G3").Value = "69:54:51"
$ws1.Range("G4").Value = "92:54:51"

# --- R2 ---
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12111:45:53"
$ws2.Range("G3").Value = "3241:29:22"
$ws2.Range("G4").Value = "479:40:56"

# --- R3: no Elapsed Duration changes ---

# --- R4 ---
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2957:35:42"
$ws4.Range("G3").Value = "184:47:57"
$ws4.Range("G4").Value = "73:00:22"
$ws4.Range("G5").Value = "70:37:55"

# --- R5 ---
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "431:34:41"

# --- R6 ---
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "72:06:59"
